# Post-study scores workbook update:
#  - relabel the "emotionSignificant" row (row 24, column A) as the csv
#    file name that actually produced it
#  - record two new experiment rows: "zoo" (emotion, no daycount) and the
#    emotion+daycount follow-up run
#  - move the active selection down to the newly added row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24 label was a placeholder ("emotionSignificant"); rename to the
# actual source file now that the daycount-validated run has its own rows.
$ws.Range("A24").Value = "everydayCountSignificant.csv"

# New row 27: "zoo" run - emotion features, no daycount
$ws.Range("A27").Value = "zoo"
$ws.Range("B27").Value = 0.97146809999999995
$ws.Range("D27").Value = 0.90649000000000002
$ws.Range("F27").Value = "emotion， 但是没有daycount"

# New row 28: emotion+daycount follow-up run (same group as row 27)
$ws.Range("B28").Value = 0.99303450000000004
$ws.Range("D28").Value = 0.90966000000000002
$ws.Range("F28").Value = "emotion+daycount"

# Column C (C27/C28) keeps the "General" look the sheet has used for every
# other row since row 19, instead of picking up the column's 0.00000 number
# format - copy the format from an already-general cell before writing the
# value so the style stays unset, same as its neighbours.
$ws.Range("C19").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 0.92774000000000001

$ws.Range("C19").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 0.93023

$ws.Application.CutCopyMode = $false

# Selection follows the freshly entered data, as it did after row 24 before.
$ws.Range("C28").Select()
